# Update release metadata on the "Metadata" sheet and swap the two
# "Mapping" columns (AK/AL) on the "Elements" sheet, per the
# release-notes.md-driven refresh of this StructureDefinition IG export.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metadata sheet: bump version/status/date, fill in the real contact.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"            # Version
$meta.Range("B6").Value  = "draft"                       # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"   # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)" # Contact

# ---------------------------------------------------------------------
# Elements sheet: the "Mapping: RIM Mapping" column and the
# "Mapping: Spécification métier vers l'extension ROR ClosingType"
# column traded places (AK <-> AL), header, data and column width alike.
# ---------------------------------------------------------------------
$els = $wb.Worksheets.Item("Elements")

# Header row (row 1)
$els.Range("AK1").Value = "Mapping: Spécification métier vers l'extension ROR ClosingType"
$els.Range("AL1").Value = "Mapping: RIM Mapping"

# Data rows
$els.Range("AK2").Value = $null
$els.Range("AL2").Value = $null

$els.Range("AK3").Value = $null
$els.Range("AL3").Value = "n/a"

$els.Range("AK4").Value = $null
$els.Range("AL4").Value = $null

$els.Range("AK5").Value = $null
$els.Range("AL5").Value = "N/A"

$els.Range("AK6").Value = "typeFermeture"
$els.Range("AL6").Value = "N/A"

# Column widths also swapped along with the columns' contents.
$els.Columns.Item(37).ColumnWidth = 67.61328125
$els.Columns.Item(38).ColumnWidth = 24.98046875
